$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new year column AN (2023) - copy formatting (bold/border/center)
# from the neighboring year-header cell AM1, then set the text value.
$ws.Range("AM1").Copy()
$ws.Range("AN1").PasteSpecial(-4122)
$ws.Range("AN1").Value = "'2023"

# New data values for 2023, rows 2-10
$ws.Range("AN2").Value = 3291
$ws.Range("AN3").Value = 2286
$ws.Range("AN4").Value = 5577
$ws.Range("AN5").Value = 1869
$ws.Range("AN6").Value = 2251
$ws.Range("AN7").Value = 4119
$ws.Range("AN8").Value = 1422
$ws.Range("AN9").Value = 35
$ws.Range("AN10").Value = 1458
